# Reorder the "Sector" rows (Onshore / Photovoltaic / Offshore wind plants)
# on every year sheet: the new row order is
#   row 5 -> Offshore wind plants
#   row 6 -> Onshore wind plants
#   row 7 -> Photovoltaic plants
# and each row's data (column E; D/F/G are always 0) travels with its label,
# which amounts to rotating the E5:E7 values down by one (E5<-old E7,
# E6<-old E5, E7<-old E6).

$wb = $excel.ActiveWorkbook

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)

    $e5 = $ws.Range("E5").Value()
    $e6 = $ws.Range("E6").Value()
    $e7 = $ws.Range("E7").Value()

    $ws.Range("C5").Value = "Offshore wind plants"
    $ws.Range("C6").Value = "Onshore wind plants"
    $ws.Range("C7").Value = "Photovoltaic plants"

    $ws.Range("E5").Value = $e7
    $ws.Range("E6").Value = $e5
    $ws.Range("E7").Value = $e6
}
